$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 29.46595573425293
$ws.Range("D2").Value = -0.08404426574706747
$ws.Range("E2").Value = 0.007063438604963698
$ws.Range("C3").Value = 29.70925712585449
$ws.Range("D3").Value = -0.04074287414550781
$ws.Range("E3").Value = 0.001659981793636689
$ws.Range("C4").Value = 29.85741233825684
$ws.Range("D4").Value = 0.01741233825683253
$ws.Range("E4").Value = 0.0003031895235703536
$ws.Range("C5").Value = 29.88316917419434
$ws.Range("D5").Value = 0.07316917419433366
$ws.Range("E5").Value = 0.005353728052280743
$ws.Range("C6").Value = 29.93117332458496
$ws.Range("D6").Value = 0.01117332458495923
$ws.Range("E6").Value = 0.0001248431822808544
$ws.Range("C7").Value = 30.08595848083496
$ws.Range("D7").Value = 0.105958480834957
$ws.Range("E7").Value = 0.01122719966085194
$ws.Range("B8").Value = 30.03999999999999
$ws.Range("C8").Value = 30.22243309020996
$ws.Range("D8").Value = 0.1824330902099689
$ws.Range("E8").Value = 0.03328183240355865
$ws.Range("B9").Value = 30.21000000000001
$ws.Range("C9").Value = 30.27715873718262
$ws.Range("D9").Value = 0.06715873718260923
$ws.Range("E9").Value = 0.004510295979962779
$ws.Range("C10").Value = 30.35295677185059
$ws.Range("D10").Value = 0.1329567718505871
$ws.Range("E10").Value = 0.01767750318092906
$ws.Range("C11").Value = 30.36369132995605
$ws.Range("D11").Value = -0.01630867004394077
$ws.Range("E11").Value = 0.0002659727186021309
$ws.Range("C12").Value = 30.49990653991699
$ws.Range("D12").Value = 0.05990653991699446
$ws.Range("E12").Value = 0.003588793524826451
$ws.Range("C13").Value = 30.43221092224121
$ws.Range("D13").Value = -0.04778907775879304
$ws.Range("E13").Value = 0.002283795953035968
$ws.Range("C14").Value = 30.43744468688965
$ws.Range("D14").Value = -0.2525553131103493
$ws.Range("E14").Value = 0.06378418618026657
$ws.Range("C15").Value = 30.48656272888184
$ws.Range("D15").Value = -0.2634372711181641
$ws.Range("E15").Value = 0.06939919581418508
$ws.Range("C16").Value = 30.66298484802246
$ws.Range("D16").Value = -0.2770151519775368
$ws.Range("E16").Value = 0.0767373944251378
$ws.Range("C17").Value = 30.77135276794434
$ws.Range("D17").Value = -0.1786472320556669
$ws.Range("E17").Value = 0.0319148335211513
$ws.Range("C18").Value = 31.13797569274902
$ws.Range("D18").Value = 0.1179756927490274
$ws.Range("E18").Value = 0.01391826407961292
$ws.Range("C19").Value = 31.23985862731934
$ws.Range("D19").Value = 0.1198586273193314
$ws.Range("E19").Value = 0.01436609054287437
$ws.Range("C20").Value = 31.34793663024902
$ws.Range("D20").Value = 0.0679366302490223
$ws.Range("E20").Value = 0.004615385729592372
$ws.Range("C21").Value = 31.21634483337402
$ws.Range("D21").Value = -0.163655166625972
$ws.Range("E21").Value = 0.02678301356337466
$ws.Range("C22").Value = 31.37014579772949
$ws.Range("D22").Value = -0.2098542022705061
$ws.Range("E22").Value = 0.04403878621059049
$ws.Range("B23").Value = 31.65000000000001
$ws.Range("C23").Value = 31.93939781188965
$ws.Range("D23").Value = 0.2893978118896428
$ws.Range("E23").Value = 0.08375109352651305
$ws.Range("C24").Value = 32.42129135131836
$ws.Range("D24").Value = 0.5412913513183639
$ws.Range("E24").Value = 0.2929963270120605
$ws.Range("C25").Value = 32.3786506652832
$ws.Range("D25").Value = 0.09865066528320199
$ws.Range("E25").Value = 0.009731953760818354
$ws.Range("C26").Value = 32.47011184692383
$ws.Range("D26").Value = 0.02011184692382528
$ws.Range("E26").Value = 0.0004044863866873805
$ws.Range("B27").Value = 32.84999999999999
$ws.Range("C27").Value = 32.71941757202148
$ws.Range("D27").Value = -0.1305824279785099
$ws.Range("E27").Value = 0.01705177049676274
$ws.Range("B28").Value = 32.90000000000001
$ws.Range("C28").Value = 32.95898818969727
$ws.Range("D28").Value = 0.05898818969725994
$ws.Range("E28").Value = 0.003479606523759924
$ws.Range("B29").Value = 33.09999999999999
$ws.Range("C29").Value = 32.89831924438477
$ws.Range("D29").Value = -0.2016807556152287
$ws.Range("E29").Value = 0.0406751271855296
$ws.Range("B30").Value = 33.40000000000001
$ws.Range("C30").Value = 33.6444206237793
$ws.Range("D30").Value = 0.2444206237792912
$ws.Range("E30").Value = 0.05974144132865781
$ws.Range("C31").Value = 33.67840194702148
$ws.Range("D31").Value = -0.02159805297851847
$ws.Range("E31").Value = 0.0004664758924628904
$ws.Range("B32").Value = 34.09999999999999
$ws.Range("C32").Value = 33.89046096801758
$ws.Range("D32").Value = -0.2095390319824162
$ws.Range("E32").Value = 0.04390660592412803
$ws.Range("B33").Value = 34.40000000000001
$ws.Range("C33").Value = 34.41019821166992
$ws.Range("D33").Value = 0.01019821166991619
$ws.Range("E33").Value = 0.0001040035212644148
$ws.Range("B34").Value = 34.90000000000001
$ws.Range("C34").Value = 35.07607650756836
$ws.Range("D34").Value = 0.1760765075683537
$ws.Range("E34").Value = 0.03100293651746851
$ws.Range("C35").Value = 35.64742660522461
$ws.Range("D35").Value = 0.3474266052246122
$ws.Range("E35").Value = 0.1207052460178985
$ws.Range("C36").Value = 35.9569206237793
$ws.Range("D36").Value = 0.256920623779294
$ws.Range("E36").Value = 0.06600820692314155
$ws.Range("C37").Value = 35.85188674926758
$ws.Range("D37").Value = -0.448113250732419
$ws.Range("E37").Value = 0.2008054854819759
$ws.Range("C38").Value = 36.47769546508789
$ws.Range("D38").Value = -0.3223045349121065
$ws.Range("E38").Value = 0.1038802132249093
$ws.Range("C39").Value = 37.24319839477539
$ws.Range("D39").Value = -0.05680160522460653
$ws.Range("E39").Value = 0.003226422356092048
$ws.Range("B40").Value = 37.90000000000001
$ws.Range("C40").Value = 38.02762222290039
$ws.Range("D40").Value = 0.1276222229003849
$ws.Range("E40").Value = 0.01628743177803554
$ws.Range("C41").Value = 38.46606826782227
$ws.Range("D41").Value = -0.03393173217773438
$ws.Range("E41").Value = 0.001151362448581494
$ws.Range("B42").Value = 38.90000000000001
$ws.Range("C42").Value = 39.12563705444336
$ws.Range("D42").Value = 0.2256370544433537
$ws.Range("E42").Value = 0.05091208033787296
$ws.Range("B43").Value = 39.40000000000001
$ws.Range("C43").Value = 39.65830612182617
$ws.Range("D43").Value = 0.2583061218261662
$ws.Range("E43").Value = 0.06672205257287421
$ws.Range("B44").Value = 39.90000000000001
$ws.Range("C44").Value = 39.70235824584961
$ws.Range("D44").Value = -0.1976417541503963
$ws.Range("E44").Value = 0.0390622629836457
$ws.Range("B45").Value = 40.09999999999999
$ws.Range("C45").Value = 39.89709854125977
$ws.Range("D45").Value = -0.2029014587402287
$ws.Range("E45").Value = 0.04116900195891272
$ws.Range("B46").Value = 40.59999999999999
$ws.Range("C46").Value = 40.38924026489258
$ws.Range("D46").Value = -0.2107597351074162
$ws.Range("E46").Value = 0.04441966594254824
$ws.Range("B47").Value = 40.90000000000001
$ws.Range("C47").Value = 40.60507583618164
$ws.Range("D47").Value = -0.2949241638183651
$ws.Range("E47").Value = 0.08698026240396182
$ws.Range("B48").Value = 41.20000000000001
$ws.Range("C48").Value = 41.06470108032227
$ws.Range("D48").Value = -0.1352989196777443
$ws.Range("E48").Value = 0.01830579766596471
$ws.Range("C49").Value = 41.33103561401367
$ws.Range("D49").Value = -0.1689643859863281
$ws.Range("E49").Value = 0.02854896373173688
$ws.Range("C50").Value = 42.09722518920898
$ws.Range("D50").Value = 0.2972251892089872
$ws.Range("E50").Value = 0.08834281310031825
$ws.Range("C51").Value = 42.46757125854492
$ws.Range("D51").Value = 0.267571258544919
$ws.Range("E51").Value = 0.0715943783993119
$ws.Range("C52").Value = 0.006692657470672714
$ws.Range("E52").Value = 2.06433120004918
$ws.Range("E53").Value = 0.04128662400098359
